$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the "password" column (E), shifting it to F.
$ws.Range("E1").EntireColumn.Insert()

# Match the new column's width to column D (role_id/nip column) width.
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# New header + value for the inserted "email" column.
$ws.Range("E1").Value = "email"
$ws.Range("E2").Value = "dosen@gmail.com"

# Turn the email value into a mailto hyperlink (also applies the built-in
# Hyperlink cell style used by Excel).
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:dosen@gmail.com")

# Match the selection left behind by the edit.
$null = $ws.Range("E7").Select()
